# Update the "Etat Virement" sheet: replace the data for rows 2-4, turn the
# former totals row (row 6) into row 5 (deleting the old, now-redundant
# row 5), and shrink the used range from A1:K6 to A1:K5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 5 entirely; this shifts the old row 6 (the blank/total row) up
# to become the new row 5, matching the target layout.
$ws.Rows.Item(5).Delete()

# Account-number column (C) holds long purely-numeric strings that must stay
# text (otherwise Excel coerces them to floating point / scientific
# notation). Force text formatting before writing those values.
$ws.Range("C2:C4").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "JALAL MED"
$ws.Range("B2").Value = "J2545456"
$ws.Range("C2").Value = "265463456416547645654645"
$ws.Range("D2").Value = "AGG"
$ws.Range("E2").Value = "BP"
$ws.Range("F2").Value = "Supervision"
$ws.Range("G2").Value = "800/SUP 1"
$ws.Range("H2").Value = "mensuelle"
$ws.Range("I2").Value = 8000
$ws.Range("J2").Value = 800
$ws.Range("K2").Value = 7200

# Row 3
$ws.Range("A3").Value = "ACHENGLI LAILA"
$ws.Range("B3").Value = "J207703"
$ws.Range("C3").Value = "00101211115087750001201090"
$ws.Range("D3").Value = "Ait souss"
$ws.Range("E3").Value = "BP Centre Sud"
$ws.Range("F3").Value = "Direction régionale"
$ws.Range("G3").Value = "901/FES "
$ws.Range("H3").Value = "mensuelle"
$ws.Range("I3").Value = 5000
$ws.Range("J3").Value = 500
$ws.Range("K3").Value = 4500

# Row 4
$ws.Range("A4").Value = "CHARIJI ABDELLAH"
$ws.Range("B4").Value = "BJ36877"
$ws.Range("C4").Value = "00101211111292695000201732"
$ws.Range("D4").Value = "AOURIR"
$ws.Range("E4").Value = "BP CENTRE SUD"
$ws.Range("F4").Value = "Logement de fonction"
$ws.Range("G4").Value = "901/LF/FES "
$ws.Range("H4").Value = "mensuelle"
$ws.Range("I4").Value = 6000
$ws.Range("J4").Value = 600
$ws.Range("K4").Value = 5400

# Row 5 (former row 6, now shifted up): keep blank text cells, update totals.
$ws.Range("A5").Value = " "
$ws.Range("B5").Value = " "
$ws.Range("C5").Value = " "
$ws.Range("D5").Value = " "
$ws.Range("E5").Value = " "
$ws.Range("F5").Value = " "
$ws.Range("G5").Value = " "
$ws.Range("H5").Value = " "
$ws.Range("I5").Value = 19000
$ws.Range("J5").Value = 1900
$ws.Range("K5").Value = 17100
